$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.136.26"
$ws.Range("E2").Value = "'  -0.03%  "
$ws.Range("D3").Value = "'1.831.92"
$ws.Range("E3").Value = "'  -0.37%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D5").Value = "'241.20"
$ws.Range("E5").Value = "'  +0.47%  "
$ws.Range("D6").Value = "'0.6635"
$ws.Range("E6").Value = "'  -2.59%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("D8").Value = "'0.07409"
$ws.Range("E8").Value = "'  -0.54%  "
$ws.Range("D9").Value = "'0.2934"
$ws.Range("E9").Value = "'  -1.85%  "
$ws.Range("D10").Value = "'22.67"
$ws.Range("E10").Value = "'  -2.37%  "
$ws.Range("D11").Value = "'0.07740"
$ws.Range("D12").Value = "'1.839.43"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("D13").Value = "'4.979"
$ws.Range("E13").Value = "'  -0.95%  "
$ws.Range("D14").Value = "'0.6671"
$ws.Range("E14").Value = "'  -2.01%  "
$ws.Range("D15").Value = "'82.63"
$ws.Range("E15").Value = "'  -5.23%  "
$ws.Range("D16").Value = "'6.081"
$ws.Range("E16").Value = "'  -1.07%  "
$ws.Range("D17").Value = "'0.000008333"
$ws.Range("E17").Value = "'  +1.65%  "
$ws.Range("D18").Value = "'29.178.18"
$ws.Range("E18").Value = "'  +0.19%  "
$ws.Range("D19").Value = "'226.85"
$ws.Range("E19").Value = "'  -1.42%  "
$ws.Range("D20").Value = "'12.46"
$ws.Range("E20").Value = "'  -0.40%  "
$ws.Range("E21").Value = "'  +0.16%  "
$ws.Range("D22").Value = "'7.148"
$ws.Range("E22").Value = "'  -2.74%  "
$ws.Range("E23").Value = "'  +0.07%  "
$ws.Range("D24").Value = "'159.37"
$ws.Range("E24").Value = "'  -1.30%  "
$ws.Range("D25").Value = "'8.599"
$ws.Range("E25").Value = "'  -1.29%  "
$ws.Range("D26").Value = "'0.1400"
$ws.Range("E26").Value = "'  -2.07%  "
$ws.Range("D27").Value = "'17.93"
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("D28").Value = "'1.507"
$ws.Range("E28").Value = "'  +0.44%  "
$ws.Range("D29").Value = "'4.105"
$ws.Range("E29").Value = "'  -3.43%  "
$ws.Range("D30").Value = "'4.027"
$ws.Range("E30").Value = "'  -2.86%  "
$ws.Range("E31").Value = "'  +0.19%  "
$ws.Range("D32").Value = "'0.05310"
$ws.Range("E32").Value = "'  -0.71%  "
$ws.Range("D33").Value = "'1.866"
$ws.Range("E33").Value = "'  +0.99%  "
$ws.Range("D34").Value = "'0.7511"
$ws.Range("E34").Value = "'  -0.27%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "'  +0.36%  "
$ws.Range("D36").Value = "'2.651"
$ws.Range("E36").Value = "'  -1.16%  "
$ws.Range("D37").Value = "'1.277.57"
$ws.Range("E37").Value = "'  -2.55%  "
$ws.Range("D38").Value = "'0.01793"
$ws.Range("E38").Value = "'  -1.60%  "
$ws.Range("D39").Value = "'2.727"
$ws.Range("E39").Value = "'  +0.32%  "
$ws.Range("D40").Value = "'0.9272"
$ws.Range("E40").Value = "'  -1.49%  "
$ws.Range("D41").Value = "'0.08560"
$ws.Range("E41").Value = "'  +3.94%  "
$ws.Range("D42").Value = "'5.915"
$ws.Range("E42").Value = "'  -2.61%  "
$ws.Range("D43").Value = "'0.9997"
$ws.Range("E43").Value = "'  +0.03%  "
$ws.Range("D44").Value = "'101.89"
$ws.Range("E44").Value = "'  -3.07%  "
$ws.Range("D45").Value = "'2.002.58"
$ws.Range("E45").Value = "'  +1.14%  "
$ws.Range("D46").Value = "'0.5148"
$ws.Range("E46").Value = "'  -0.53%  "
$ws.Range("D47").Value = "'1.761"
$ws.Range("E47").Value = "'  -0.77%  "
$ws.Range("E48").Value = "'  -0.99%  "
$ws.Range("D49").Value = "'63.26"
$ws.Range("E49").Value = "'  -1.53%  "
$ws.Range("D50").Value = "'0.05889"
$ws.Range("E50").Value = "'  -0.87%  "
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.820"
$ws.Range("E51").Value = "'  -6.47%  "
